$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (row 1) ---
# Row 1 used to (incorrectly) duplicate row 2's data; replace it with real headers,
# matching the other property-type sheets (name/capacity/owner/register_date/...).
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (HYUNDAI car) ---
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2012-02-10"
$ws.Cells.Item(2, 11).Value = "盧秀燕"
$ws.Cells.Item(2, 12).Value = 869
$ws.Cells.Item(2, 13).Value = "tmp61a71"
$ws.Cells.Item(2, 14).Value = 32

# --- Row 3 (國瑞 car) ---
$ws.Cells.Item(3, 2).Value = "國瑞"
$ws.Cells.Item(3, 5).Value = "99年08月16日"
$ws.Cells.Item(3, 8).Value = "land"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 10).Value = "2012-02-10"
$ws.Cells.Item(3, 11).Value = "盧秀燕"
$ws.Cells.Item(3, 12).Value = 869
$ws.Cells.Item(3, 13).Value = "tmp61a71"
$ws.Cells.Item(3, 14).Value = 33
